# ESCALETA FINAL CN_11_09_CO.xlsx - "Cambio de F a la izquierda"
# Swap the "M" (Masculino) and "F" (Femenino) engine-type columns (M and N)
# so that the F-list now lives in column M (to the left of M-list, now in N).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------
# 1) Swap the cell VALUES of columns M and N for every data row (2-92).
#    Row 1 is a merged M1:N1 header ("Tipo de Motor") and must stay put.
# ---------------------------------------------------------------------
$mRange = $ws.Range("M2:M92")
$nRange = $ws.Range("N2:N92")
$mVals = $mRange.Value()
$nVals = $nRange.Value()
$mRange.Value = $nVals
$nRange.Value = $mVals

# ---------------------------------------------------------------------
# 2) Row 13 carries a slightly different (non-wrap) style on whichever
#    cell holds the value; keep that nuance attached to the data as it
#    moves from M13 to N13.
# ---------------------------------------------------------------------
$ws.Range("M13").WrapText = $true
$ws.Range("N13").WrapText = $false

# ---------------------------------------------------------------------
# 3) Swap the column widths (and the "extra" style carried past the
#    used range) between M and N so formatting follows the data.
# ---------------------------------------------------------------------
$colM = $ws.Columns.Item(13)
$colN = $ws.Columns.Item(14)
$mWidth = $colM.ColumnWidth()
$nWidth = $colN.ColumnWidth()
$colM.ColumnWidth = $nWidth
$colN.ColumnWidth = $mWidth

# ---------------------------------------------------------------------
# 4) Swap which column uses the list validated against "F" (DATOS!B)
#    and which uses "M" (DATOS!C). M should now validate against the
#    F list, N against the M list.
# ---------------------------------------------------------------------
$ws.Range("M3:M92").Validation.Add(3, 1, 1, "=DATOS!`$B`$2:`$B`$14")
$ws.Range("N3:N92").Validation.Add(3, 1, 1, "=DATOS!`$C`$2:`$C`$39")

# ---------------------------------------------------------------------
# 5) Turn off the AutoFilter (the filter buttons/range are removed).
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false

# ---------------------------------------------------------------------
# 6) Update the window/view state: scroll so column J is at the left
#    edge, and move the active selection to O11.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("O11").Select() | Out-Null
